$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 111 ("Vega Modelo de
# Temuco" / Albahaca weekly records), shifting all existing rows (111-175)
# down to (113-177). Excel's Rows.Insert carries the row-above formatting
# along (keeps the date number format on column D), matching the target
# dimension A1:R177.
$ws.Rows.Item(111).Insert()
$ws.Rows.Item(111).Insert()

# Populate the two newly inserted rows with the new weekly price records.
# Columns A,B,C,E,F,G,H,I,R are constant for every record in this sheet.

# Row 111: Región de La Araucanía
$ws.Cells.Item(111, 1).Value = 10
$ws.Cells.Item(111, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(111, 3).Value = "La Araucanía"
$ws.Cells.Item(111, 4).Value = 44574
$ws.Cells.Item(111, 5).Value = 9
$ws.Cells.Item(111, 6).Value = 100112052
$ws.Cells.Item(111, 7).Value = "Albahaca"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 110
$ws.Cells.Item(111, 11).Value = 5000
$ws.Cells.Item(111, 12).Value = 5000
$ws.Cells.Item(111, 13).Value = 5000
$ws.Cells.Item(111, 14).Value = "$/paquete"
$ws.Cells.Item(111, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(111, 16).Value = 5000
$ws.Cells.Item(111, 17).Value = 1
$ws.Cells.Item(111, 18).Value = "Hortaliza"

# Row 112: Región del Maule
$ws.Cells.Item(112, 1).Value = 10
$ws.Cells.Item(112, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(112, 3).Value = "La Araucanía"
$ws.Cells.Item(112, 4).Value = 44574
$ws.Cells.Item(112, 5).Value = 9
$ws.Cells.Item(112, 6).Value = 100112052
$ws.Cells.Item(112, 7).Value = "Albahaca"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 55
$ws.Cells.Item(112, 11).Value = 5000
$ws.Cells.Item(112, 12).Value = 5000
$ws.Cells.Item(112, 13).Value = 5000
$ws.Cells.Item(112, 14).Value = "$/paquete"
$ws.Cells.Item(112, 15).Value = "Región del Maule"
$ws.Cells.Item(112, 16).Value = 5000
$ws.Cells.Item(112, 17).Value = 1
$ws.Cells.Item(112, 18).Value = "Hortaliza"
